# Auto-generated edit script: add 23 new localization rows (59-81) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "popSave_01emptyNodes"
$ws.Range("B59").Value = "popup save button (popSave)"
$ws.Range("C59").Value = "Please return to your Cognitive-Affective Map and add text to the empty concepts."
$ws.Range("D59").Value = "Bitte kehren Sie zu Ihrer Cognitive-Affective Map zurück und ergänzen sie Text zu den leeren Konzepten."
$ws.Rows.Item(59).RowHeight = 75

$ws.Range("A60").Value = "popSave_02emptyNodes"
$ws.Range("B60").Value = "popup save button (popSave)"
$ws.Range("C60").Value = " concept(s) are empty."
$ws.Range("D60").Value = " Konzept(e) sind leer."
$ws.Rows.Item(60).RowHeight = 30

$ws.Range("A61").Value = "popSave_01numNodes"
$ws.Range("B61").Value = "popup save button (popSave)"
$ws.Range("C61").Value = "Please return to your Cognitive-Affective Map and add additional concepts to it."
$ws.Range("D61").Value = "Bitte kehren Sie zu Ihrer Cognitive Affective Map zurück und fügen Sie weitere Konzepte hinzu."
$ws.Rows.Item(61).RowHeight = 60

$ws.Range("A62").Value = "popSave_02numNodes"
$ws.Range("B62").Value = "popup save button (popSave)"
$ws.Range("C62").Value = "Please draw at least "
$ws.Range("D62").Value = "Bitte zeichnen Sie mindestens "
$ws.Rows.Item(62).RowHeight = 30

$ws.Range("A63").Value = "popSave_03numNodes"
$ws.Range("B63").Value = "popup save button (popSave)"
$ws.Range("C63").Value = " concepts."
$ws.Range("D63").Value = " Konzepte."
$ws.Rows.Item(63).RowHeight = 30

$ws.Range("A64").Value = "popSave_01unconnectedA"
$ws.Range("B64").Value = "popup save button (popSave)"
$ws.Range("C64").Value = "Please return to your Cognitive-Affective Map and add additional connections to it."
$ws.Range("D64").Value = "Bitte kehren Sie zu Ihrer Cognitive-Affective Map zurück und fügen Sie weitere Verbindungen hinzu."
$ws.Rows.Item(64).RowHeight = 75

$ws.Range("A65").Value = "popSave_02unconnectedA"
$ws.Range("B65").Value = "popup save button (popSave)"
$ws.Range("C65").Value = "Please connect all your concepts within your Cognitive-Affective Map."
$ws.Range("D65").Value = "Bitte verbinden Sie alle Ihre Konzepte innerhalb Ihrer Cognitive-Affective Map."
$ws.Rows.Item(65).RowHeight = 60

$ws.Range("A66").Value = "popSave_01unconnectedB"
$ws.Range("B66").Value = "popup save button (popSave)"
$ws.Range("C66").Value = "Please return to your Cognitive-Affective Map and add additional connections to it."
$ws.Range("D66").Value = "Bitte kehren Sie zu Ihrer Cognitive-Affective Map zurück und fügen Sie weitere Verbindungen hinzu."
$ws.Rows.Item(66).RowHeight = 75

$ws.Range("A67").Value = "popSave_02unconnectedB"
$ws.Range("B67").Value = "popup save button (popSave)"
$ws.Range("C67").Value = "Please connect all your "
$ws.Range("D67").Value = "Bitte verbinden Sie alle Ihre"
$ws.Rows.Item(67).RowHeight = 30

$ws.Range("A68").Value = "popSave_03unconnectedB"
$ws.Range("B68").Value = "popup save button (popSave)"
$ws.Range("C68").Value = " distinct groups of concepts within your Cognitive-Affective Map."
$ws.Range("D68").Value = " verschiedene Gruppen von Konzepten innerhalb Ihrer Cognitive-Affective Map."
$ws.Rows.Item(68).RowHeight = 45

$ws.Range("A69").Value = "popSave_01savedData"
$ws.Range("B69").Value = "popup save button (popSave)"
$ws.Range("C69").Value = "Your CAM data has been sent to the sever. Thank you for participating! You will be forwarded to the end or the final part of the study."
$ws.Range("D69").Value = "Ihre CAM-Daten wurden an den Sever gesendet. Wir danken Ihnen für Ihre Teilnahme! Sie werden zum Ende oder zum letzten Teil der Studie weitergeleitet."
$ws.Rows.Item(69).RowHeight = 105

$ws.Range("A70").Value = "popSave_01notSavedData"
$ws.Range("B70").Value = "popup save button (popSave)"
$ws.Range("C70").Value = "You would have send the CAM data successfully to a sever. To save permanently your data please use our administrative panel or host the C.A.M.E.L. software on your own server."
$ws.Range("D70").Value = "Sie hätten die CAM-Daten erfolgreich an einen Server gesendet. Um Ihre Daten dauerhaft zu speichern, verwenden Sie bitte unser Administrationspanel oder hosten Sie die C.A.M.E.L. Software auf Ihrem eigenen Server."
$ws.Rows.Item(70).RowHeight = 135

$ws.Range("A71").Value = "ndw_01tooManyWords"
$ws.Range("B71").Value = "nodes dialog warnings (ndw)"
$ws.Range("C71").Value = "Instead, please draw several connected concepts."
$ws.Range("D71").Value = "Zeichnen Sie stattdessen bitte mehrere zusammenhängende Konzepte."
$ws.Rows.Item(71).RowHeight = 60

$ws.Range("A72").Value = "ndw_02tooManyWords"
$ws.Range("B72").Value = "nodes dialog warnings (ndw)"
$ws.Range("C72").Value = "Please do not use more than "
$ws.Range("D72").Value = "Bitte verwenden Sie nicht mehr als "
$ws.Rows.Item(72).RowHeight = 30

$ws.Range("A73").Value = "ndw_03tooManyWords"
$ws.Range("B73").Value = "nodes dialog warnings (ndw)"
$ws.Range("C73").Value = " words for a single concept!"
$ws.Range("D73").Value = " Wörter für ein einzelnes Konzept!"
$ws.Rows.Item(73).RowHeight = 30

$ws.Range("A74").Value = "ndw_03tooManyWordsA"
$ws.Range("B74").Value = "nodes dialog warnings (ndw)"
$ws.Range("C74").Value = " characters for a single concept!"
$ws.Range("D74").Value = " Zeichen für ein einzelnes Konzept!"
$ws.Rows.Item(74).RowHeight = 30

$ws.Range("A75").Value = "ndw_01predefinedConcept"
$ws.Range("B75").Value = "nodes dialog warnings (ndw)"
$ws.Range("C75").Value = "Instead, please choose other concepts."
$ws.Range("D75").Value = "Wählen Sie stattdessen bitte andere Konzepte."
$ws.Rows.Item(75).RowHeight = 30

$ws.Range("A76").Value = "ndw_02predefinedConcept"
$ws.Range("B76").Value = "nodes dialog warnings (ndw)"
$ws.Range("C76").Value = "You cannot change the content of predefined concepts."
$ws.Range("D76").Value = "Sie können den Inhalt der vordefinierten Konzepte nicht ändern."
$ws.Rows.Item(76).RowHeight = 45

$ws.Range("A77").Value = "ndw_01ambivalentConcept"
$ws.Range("B77").Value = "nodes dialog warnings (ndw)"
$ws.Range("C77").Value = "To change the ambivalent concept again, please uncheck the box."
$ws.Range("D77").Value = "Um den ambivalenten Begriff wieder zu ändern, deaktivieren Sie bitte das Kontrollkästchen."
$ws.Rows.Item(77).RowHeight = 60

$ws.Range("A78").Value = "edw_01notDeleteNode"
$ws.Range("B78").Value = "Elements dialog warnings (edw)"
$ws.Range("C78").Value = "Instead, please choose other concepts."
$ws.Range("D78").Value = "Bitte wählen Sie stattdessen andere Konzepte."
$ws.Rows.Item(78).RowHeight = 45

$ws.Range("A79").Value = "edw_02notDeleteNode"
$ws.Range("B79").Value = "Elements dialog warnings (edw)"
$ws.Range("C79").Value = "You cannot delete predefined concepts."
$ws.Range("D79").Value = "Sie können keine vordefinierten Konzepte löschen."
$ws.Rows.Item(79).RowHeight = 45

$ws.Range("A80").Value = "edw_01notDeleteConnector"
$ws.Range("B80").Value = "Elements dialog warnings (edw)"
$ws.Range("C80").Value = "Instead, please choose other connectors."
$ws.Range("D80").Value = "Bitte wählen Sie stattdessen andere Verbindungen."
$ws.Rows.Item(80).RowHeight = 45

$ws.Range("A81").Value = "edw_02notDeleteConnector"
$ws.Range("B81").Value = "Elements dialog warnings (edw)"
$ws.Range("C81").Value = "You cannot delete predefined connectors."
$ws.Range("D81").Value = "Sie können keine vordefinierten Verbindungen löschen."
$ws.Rows.Item(81).RowHeight = 45

# Widen column A to fit the new, longer identifier strings
$ws.Columns.Item(1).ColumnWidth = 28.333333333333332

# Restore view/selection focus near the newly added rows
$ws.Range("A79").Select()

